$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellAddr, $val) {
    # Writes $val into the cell as literal text. Cells whose content
    # looks like a plain number (e.g. "408.00", "0.999") would otherwise
    # be silently coerced to a numeric cell by Excel, which would lose
    # the original formatting (trailing zeros, exact decimal text).
    # Forcing a Text number format guarantees the literal string is kept,
    # then the format/style is reverted to the sheet default (General /
    # "Normal") afterwards so no stray cell styling is left behind.
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.NumberFormat = "General"
    $rng.Style = "Normal"
}

$ws.Range("D2").Value = '62.524.09'
$ws.Range("E2").Value = '  -0.51%  '

$ws.Range("D3").Value = '3.443.36'
$ws.Range("E3").Value = '  -1.19%  '

$ws.Range("E4").Value = '  +0.13%  '

Set-TextValue "D5" '408.00'
$ws.Range("E5").Value = '  -0.43%  '

Set-TextValue "D6" '134.14'
$ws.Range("E6").Value = '  +1.21%  '

$ws.Range("E7").Value = '  -1.87%  '

Set-TextValue "D8" '0.999'
$ws.Range("E8").Value = '  -0.01%  '

$ws.Range("E9").Value = '  -1.57%  '

$ws.Range("E10").Value = '  -4.81%  '

Set-TextValue "D11" '42.20'
$ws.Range("E11").Value = '  -2.55%  '

$ws.Range("E12").Value = '  -0.77%  '

$ws.Range("E13").Value = '  -3.61%  '

Set-TextValue "D14" '19.95'
$ws.Range("E14").Value = '  -1.88%  '

$ws.Range("D15").Value = '3.468.97'
$ws.Range("E15").Value = '  +0.00%  '

$ws.Range("D16").Value = '62.450.05'
$ws.Range("E16").Value = '  -0.38%  '

Set-TextValue "D17" '11.39'
$ws.Range("E17").Value = '  +4.02%  '

$ws.Range("E18").Value = '  -2.87%  '

Set-TextValue "D19" '0.0000134'
$ws.Range("E19").Value = '  -3.53%  '

$ws.Range("E20").Value = '  -5.76%  '

Set-TextValue "D21" '84.16'
$ws.Range("E21").Value = '  +1.45%  '

Set-TextValue "D22" '314.85'
$ws.Range("E22").Value = '  +0.65%  '

$ws.Range("E23").Value = '  -2.22%  '

$ws.Range("E24").Value = '  -1.09%  '

Set-TextValue "D25" '4.73'
$ws.Range("E25").Value = '  +7.99%  '

$ws.Range("E26").Value = '  -2.55%  '

$ws.Range("E27").Value = '  -0.98%  '

Set-TextValue "D28" '2.79'
$ws.Range("E28").Value = '  +4.68%  '

$ws.Range("E29").Value = '  -3.01%  '

Set-TextValue "D30" '0.173'
$ws.Range("E30").Value = '  -3.81%  '

Set-TextValue "D31" '0.116'
$ws.Range("E31").Value = '  -4.44%  '

Set-TextValue "D32" '42.37'
$ws.Range("E32").Value = '  -2.37%  '

$ws.Range("E33").Value = '  -0.10%  '

$ws.Range("E34").Value = '  -4.90%  '

$ws.Range("E35").Value = '  -2.08%  '

Set-TextValue "D36" '51.49'
$ws.Range("E36").Value = '  -2.06%  '

Set-TextValue "D37" '0.998'
$ws.Range("E37").Value = '  +0.19%  '

$ws.Range("E38").Value = '  -5.92%  '

$ws.Range("B39").Value = 'TheGraph'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
Set-TextValue "D39" '0.323'
$ws.Range("E39").Value = '  +12.16%  '

$ws.Range("B40").Value = 'Stacks'
$ws.Range("C40").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
Set-TextValue "D40" '2.96'
$ws.Range("E40").Value = '  -2.20%  '

$ws.Range("B41").Value = 'Monero'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue "D41" '138.23'
$ws.Range("E41").Value = '  -0.04%  '

$ws.Range("B42").Value = 'ARBITRUM'
$ws.Range("C42").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue "D42" '1.99'
$ws.Range("E42").Value = '  -0.59%  '

$ws.Range("E43").Value = '  -1.06%  '

Set-TextValue "D44" '4.02'
$ws.Range("E44").Value = '  -0.10%  '

Set-TextValue "D45" '16.85'
$ws.Range("E45").Value = '  -4.90%  '

Set-TextValue "D46" '2.23'
$ws.Range("E46").Value = '  -1.12%  '

Set-TextValue "D47" '21.39'
$ws.Range("E47").Value = '  -5.11%  '

$ws.Range("D48").Value = '2.129.00'
$ws.Range("E48").Value = '  -3.72%  '

$ws.Range("E49").Value = '  -3.41%  '

$ws.Range("E50").Value = '  +20.15%  '

$ws.Range("E51").Value = '  +1.64%  '
